$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.715.46'
$ws.Range("E2").Value = '  -2.55%  '

$ws.Range("D3").Value = '1.559.41'
$ws.Range("E3").Value = '  -0.19%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").Value = '205.82'
$ws.Range("E5").Value = '  -1.04%  '

$ws.Range("E6").Value = '  -2.06%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  +0.78%  '

$ws.Range("E9").Value = '  -0.39%  '

$ws.Range("E10").Value = '  -1.01%  '

$ws.Range("E11").Value = '  -0.33%  '

$ws.Range("D12").Value = '1.781.85'
$ws.Range("E12").Value = '  -0.19%  '

$ws.Range("D13").Value = '1.563.71'
$ws.Range("E13").Value = '  -0.23%  '

$ws.Range("E14").Value = '  -2.03%  '

$ws.Range("E15").Value = '  -0.72%  '

$ws.Range("D16").Value = '61.43'
$ws.Range("E16").Value = '  -2.89%  '

$ws.Range("D17").Value = '26.743.04'
$ws.Range("E17").Value = '  -2.47%  '

$ws.Range("D18").Value = '214.31'
$ws.Range("E18").Value = '  +0.89%  '

$ws.Range("D19").Value = '7.33'
$ws.Range("E19").Value = '  +1.39%  '

$ws.Range("E20").Value = '  -1.64%  '

$ws.Range("E21").Value = '  +0.05%  '

$ws.Range("E22").Value = '  -0.36%  '

$ws.Range("D23").Value = '9.35'
$ws.Range("E23").Value = '  -1.64%  '

$ws.Range("D24").Value = '2.01'
$ws.Range("E24").Value = '  -0.20%  '

$ws.Range("D25").Value = '152.82'
$ws.Range("E25").Value = '  +0.02%  '

$ws.Range("D26").Value = '6.77'
$ws.Range("E26").Value = '  +0.56%  '

$ws.Range("D27").Value = '14.84'
$ws.Range("E27").Value = '  -0.86%  '

$ws.Range("E28").Value = '  +0.10%  '

$ws.Range("E29").Value = '  -1.25%  '

$ws.Range("E30").Value = '  -3.71%  '

$ws.Range("E31").Value = '  -1.46%  '

$ws.Range("E32").Value = '  -1.13%  '

$ws.Range("D33").Value = '1.387.81'
$ws.Range("E33").Value = '  +2.08%  '

$ws.Range("E34").Value = '  -1.21%  '

$ws.Range("E35").Value = '  +2.10%  '

$ws.Range("E36").Value = '  -0.51%  '

$ws.Range("D37").Value = '0.931'
$ws.Range("E37").Value = '  -4.34%  '

$ws.Range("E38").Value = '  -2.53%  '

$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = '0.516'
$ws.Range("E39").Value = '  -2.89%  '

$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").Value = '0.809'
$ws.Range("E40").Value = '  -1.25%  '

$ws.Range("D42").Value = '0.995'
$ws.Range("E42").Value = '  +2.16%  '

$ws.Range("E43").Value = '  +2.73%  '

$ws.Range("E44").Value = '  +1.77%  '

$ws.Range("D45").Value = '1.77'

$ws.Range("D46").Value = '63.20'
$ws.Range("E46").Value = '  -1.26%  '

$ws.Range("D47").Value = '1.695.04'
$ws.Range("E47").Value = '  -0.25%  '

$ws.Range("D48").Value = '85.58'
$ws.Range("E48").Value = '  +0.17%  '

$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₇0966'
$ws.Range("E49").Value = '  -1.79%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.0493'
$ws.Range("E50").Value = '  -0.11%  '

$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").Value = '0.0945'
$ws.Range("E51").Value = '  -0.85%  '

